# edit.ps1 - apply the commit's changes to the VRP / WOC / GA research paper
$d = $word.ActiveDocument

function FindReplace($findText, $replaceText, [bool]$matchCase = $true) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $matchCase, $false, $false, $false, $false, `
                       $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Title: "Vehicle Routing Problem: Using Wisdom of Artificial Crowds and Genetic Algorithms"
#    -> "Approximating solutions to the vehicle routing problem using wisdom of artificial
#        crowds with genetic algorithms"
FindReplace "Vehicle Routing Problem: Using Wisdom of Artificial Crowds and Genetic Algorithms" `
            "Approximating solutions to the vehicle routing problem using wisdom of artificial crowds with genetic algorithms"

# 2. Date: October -> November
FindReplace "(October 2019)" "(November 2019)"

# 3. Abstract: "...using Wisdom of Artificial Crowds and genetic algorithms as well as..."
#    "Crowds and " -> "Crowds with "
FindReplace "Wisdom of Artificial Crowds and genetic algorithms" `
            "Wisdom of Artificial Crowds with genetic algorithms"

# 4. Abstract closing sentence replaced
FindReplace "This algorithm was implemented in Python and tested on several datasets showing superiority over similar methods." `
            "The algorithm presented in this paper was implemented in Python and tested on several datasets producing approximations superior to any of the genetic algorithms in the crowd at the cost of post processing overhead."

Write-Output "basic text replace done"

# 5. Relocate the hidden "_GoBack" bookmark from the "Abdoun, et al." paragraph
#    to just before "post processing overhead" in the abstract (re-adding a
#    bookmark with this reserved name moves it and drops the old location).
$rngGoBack = $d.Content
$rngGoBack.Find.Execute("post processing overhead", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$goBackPoint = $d.Range($rngGoBack.Start, $rngGoBack.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

Write-Output "bookmark relocated"

# 6. Remove 3 of the 4 blank paragraphs following the abstract (keep the first).
$pBlank1 = $d.Paragraphs(7)
$pBlank3 = $d.Paragraphs(9)
$blankRange = $d.Range($pBlank1.Range.Start, $pBlank3.Range.End)
$blankRange.Delete() | Out-Null

Write-Output "blank paragraphs trimmed"

# 7. Introduction: remove the proofing mark around "Ramser" (merge runs)
FindReplace "George Dantzig and John Ramser in" "George Dantzig and John Ramser in"

# 8. Merge the "Abdoun, et al...." runs that used to carry the _GoBack bookmark
FindReplace "Abdoun, et al. uses a sequence of genetic operators, first applying a crossover method" `
            "Abdoun, et al. uses a sequence of genetic operators, first applying a crossover method"

# 9. Remove the proofing mark around "Surowiecki" (merge runs)
FindReplace "was first coined by James Surowiecki in 2004" "was first coined by James Surowiecki in 2004"

Write-Output "proofing marks merged"
